$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 19

# Row 3
$ws.Range("I3").Value = 3.8
$ws.Range("L3").Value = 4.5
$ws.Range("AC3").Value = 5.5
$ws.Range("AX3").Value = 23

# Row 4
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5

# Row 5
$ws.Range("L5").Value = 5.5
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("AC5").Value = 6.5
$ws.Range("AI5").Value = 23
$ws.Range("AM5").Value = 51
$ws.Range("AS5").Value = 251
$ws.Range("AX5").Value = 29
$ws.Range("BA5").Value = 151

# Row 6
$ws.Range("G6").Value = 3.1
$ws.Range("H6").Value = 2.7
$ws.Range("I6").Value = 2.75
$ws.Range("K6").Value = 1.8
$ws.Range("L6").Value = 3.75
$ws.Range("O6").Value = 1.67
$ws.Range("P6").Value = 2.1
$ws.Range("Q6").Value = 3.4
$ws.Range("R6").Value = 1.33
$ws.Range("S6").Value = 1.73
$ws.Range("T6").Value = 2.08
$ws.Range("U6").Value = 2.5
$ws.Range("V6").Value = 1.5
$ws.Range("W6").Value = 6
$ws.Range("X6").Value = 12
$ws.Range("AD6").Value = 6
$ws.Range("AE6").Value = 23
$ws.Range("AF6").Value = 101
$ws.Range("AH6").Value = 5.5
$ws.Range("AL6").Value = 34
$ws.Range("AO6").Value = 21
$ws.Range("AP6").Value = 41
$ws.Range("AT6").Value = 2
$ws.Range("AV6").Value = 101
$ws.Range("AX6").Value = 19
$ws.Range("AY6").Value = 41
$ws.Range("AZ6").Value = 67
$ws.Range("BA6").Value = 126

# Row 7
$ws.Range("G7").Value = 1.75
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 5
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 6.5
$ws.Range("Q7").Value = 2.5
$ws.Range("R7").Value = 1.5
$ws.Range("AD7").Value = 7

# Row 8
$ws.Range("H8").Value = 3.75
$ws.Range("N8").Value = 7.5
$ws.Range("Z8").Value = 10
$ws.Range("AA8").Value = 15
$ws.Range("AC8").Value = 7.5
$ws.Range("AI8").Value = 34
$ws.Range("AO8").Value = 8
$ws.Range("AQ8").Value = 26

# Row 9
$ws.Range("G9").Value = 1.75
$ws.Range("H9").Value = 4.1
$ws.Range("I9").Value = 4.1
$ws.Range("J9").Value = 2.25
$ws.Range("K9").Value = 2.4
$ws.Range("L9").Value = 4.33
$ws.Range("O9").Value = 1.17
$ws.Range("P9").Value = 5
$ws.Range("Q9").Value = 1.57
$ws.Range("R9").Value = 2.35
$ws.Range("AA9").Value = 13
$ws.Range("AB9").Value = 21
$ws.Range("AD9").Value = 8
$ws.Range("AI9").Value = 23
$ws.Range("AJ9").Value = 13
$ws.Range("AK9").Value = 41
$ws.Range("AQ9").Value = 26
$ws.Range("BA9").Value = 67
